$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.373.75'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.971.88'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '199.55'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '627.83'
$ws.Range("E6").Value = '  +5.30%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("D10").Value = '2.971.30'
$ws.Range("E10").Value = '  +2.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.429'
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("E13").Value = '  +1.89%  '
$ws.Range("D14").Value = '3.513.87'
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.92'
$ws.Range("E15").Value = '  +6.15%  '
$ws.Range("D16").Value = '76.305.45'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '2.965.34'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.37'
$ws.Range("E19").Value = '  +5.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.78'
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '370.96'
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("E23").Value = '  -2.07%  '
$ws.Range("E24").Value = '  +2.15%  '
$ws.Range("D25").Value = '3.122.43'
$ws.Range("E25").Value = '  +2.18%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  +2.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.68'
$ws.Range("E28").Value = '  +0.79%  '
$ws.Range("E29").Value = '  -2.18%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.21'
$ws.Range("E31").Value = '  +6.87%  '
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '505.40'
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("E34").Value = '  +7.34%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.27'
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.38'
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("E38").Value = '  +1.48%  '
$ws.Range("E39").Value = '  +10.77%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '184.85'
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.104'
$ws.Range("E41").Value = '  +13.68%  '
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.58'
$ws.Range("E44").Value = '  +6.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.91'
$ws.Range("E45").Value = '  -1.50%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.63'
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("B47").Value = 'ImmutableX'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.22'
$ws.Range("E47").Value = '  +1.79%  '
$ws.Range("E48").Value = '  +5.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.583'
$ws.Range("E49").Value = '  +1.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  +2.90%  '